$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '54.372.97'
$ws.Range('E2').Value = '  -7.28%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.871.64'
$ws.Range('E3').Value = '  -10.05%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.17%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '471.62'
$ws.Range('E5').Value = '  -11.36%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '125.90'
$ws.Range('E6').Value = '  -6.37%  '

$ws.Range('E7').Value = '  -0.08%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '2.867.50'
$ws.Range('E8').Value = '  -10.19%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.402'
$ws.Range('E9').Value = '  -11.70%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '6.63'
$ws.Range('E10').Value = '  -11.07%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0960'
$ws.Range('E11').Value = '  -14.78%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.331'
$ws.Range('E12').Value = '  -15.61%  '

$ws.Range('E13').Value = '  -4.45%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '3.351.43'
$ws.Range('E14').Value = '  -10.57%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '23.26'
$ws.Range('E15').Value = '  -9.36%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '54.252.83'
$ws.Range('E16').Value = '  -7.63%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.860.88'
$ws.Range('E17').Value = '  -10.65%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.0000134'
$ws.Range('E18').Value = '  -13.86%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '5.35'
$ws.Range('E19').Value = '  -9.16%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '11.45'
$ws.Range('E20').Value = '  -13.00%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '7.05'
$ws.Range('E21').Value = '  -13.10%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '295.63'
$ws.Range('E22').Value = '  -17.39%  '

$ws.Range('E23').Value = '  +0.07%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.444'
$ws.Range('E24').Value = '  -13.62%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '58.85'
$ws.Range('E25').Value = '  -15.61%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.993'
$ws.Range('E26').Value = '  -0.66%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.153'
$ws.Range('E27').Value = '  -9.93%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.998'
$ws.Range('E28').Value = '  -0.08%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.0₃0809'
$ws.Range('E29').Value = '  -14.68%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '6.19'
$ws.Range('E30').Value = '  -11.92%  '

$ws.Range('B31').Value = 'Fetch.AI'
$ws.Range('C31').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.14'
$ws.Range('E31').Value = '  -4.05%  '

$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '6.22'
$ws.Range('E32').Value = '  -10.82%  '

$ws.Range('E33').Value = '  -15.50%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '18.87'
$ws.Range('E34').Value = '  -12.56%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.22'
$ws.Range('E35').Value = '  -13.24%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '135.42'
$ws.Range('E36').Value = '  -16.07%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '5.43'
$ws.Range('E37').Value = '  -13.99%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.22'
$ws.Range('E38').Value = '  -13.70%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '23.12'
$ws.Range('E39').Value = '  -10.40%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0617'
$ws.Range('E40').Value = '  -11.87%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.885.48'
$ws.Range('E41').Value = '  -10.60%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.997'
$ws.Range('E42').Value = '  -0.45%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '35.02'
$ws.Range('E43').Value = '  -14.21%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.952'
$ws.Range('E44').Value = '  -13.03%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.602'
$ws.Range('E45').Value = '  -15.20%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.32'
$ws.Range('E46').Value = '  -10.73%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.39'
$ws.Range('E47').Value = '  -14.81%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.045.35'
$ws.Range('E48').Value = '  -10.48%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '5.34'
$ws.Range('E49').Value = '  -14.09%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '17.93'
$ws.Range('E50').Value = '  -12.13%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0212'
$ws.Range('E51').Value = '  -11.04%  '
